$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink the saved window height (cosmetic workbook-view state).
$excel.ActiveWindow.Height = 9000

# Re-seat the AutoFilter over the still-13-data-row range (header + rows 2..14)
# before inserting the new row, so the stored ref stays "A1:F14" rather than
# snapping to the post-insert used range.
$ws.AutoFilterMode = $false
$ws.Range("A1:F14").AutoFilter()

# Insert a new row 15 by copying row 14 (same formatting/styles) and shifting
# down, then overwrite the two text cells with the new test-case data.
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(15).Insert(-4121)  # xlShiftDown
$ws.Cells.Item(15, 1).Value = "Product_Summary-Hide/Show_account_on_Product_List_[WEB]_1"
$ws.Cells.Item(15, 2).Value = "C70773"

# Update the hidden _xlnm._FilterDatabase defined name to match the new
# filter range.
$n = $wb.Names.Item(1)
$n.RefersTo = "=Sheet1!`$A`$1:`$F`$14"

# Move the active selection to B17 (matches the post-edit workbook state).
$ws.Range("B17").Select()
